# Daily attendance processing - 2026-01-04 14:03:33
# For every already-recorded session row, the "Recorded By" (column G)
# value listing both the proctor email and "System" had the two names
# swapped from "email, System" to "System, email".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
